# Update the two-digit division problems in the worksheet table to the
# regenerated values. Several of the original expressions (e.g. "23÷4=",
# "85÷3=", "73÷4=") repeat at multiple positions in the table, so a plain
# document-wide Find/Replace would clobber the wrong cell. Instead each
# cell is targeted individually by (row, column) and only the visible
# text portion of that cell's range (i.e. excluding the trailing
# paragraph/cell-mark) is rewritten, which leaves the run/paragraph
# formatting untouched.

$d = $word.ActiveDocument
$t = $d.Tables(1)

function Set-CellText($table, $row, $col, $newText) {
    $cellRange = $table.Cell($row, $col).Range
    # Cell range includes a trailing paragraph+cell mark counted as a
    # single addressable unit; drop it so only the literal text is
    # replaced.
    $textRange = $d.Range($cellRange.Start, $cellRange.End - 1)
    $textRange.Text = $newText
}

Set-CellText $t 1 1 "97÷3="   # was 23÷4=
Set-CellText $t 1 2 "39÷3="   # was 60÷2=
Set-CellText $t 1 3 "10÷2="   # was 85÷3=
Set-CellText $t 1 4 "73÷4="   # was 77÷7=
Set-CellText $t 1 5 "68÷3="   # was 34÷9=

Set-CellText $t 5 1 "38÷6="   # was 79÷2=
Set-CellText $t 5 2 "24÷3="   # was 35÷5=
Set-CellText $t 5 3 "57÷4="   # was 85÷8=
Set-CellText $t 5 4 "51÷7="   # was 23÷4=
Set-CellText $t 5 5 "76÷3="   # was 69÷7=

Set-CellText $t 9 1 "38÷7="   # was 85÷3=
Set-CellText $t 9 2 "89÷2="   # was 68÷4=
Set-CellText $t 9 3 "47÷3="   # was 41÷3=
Set-CellText $t 9 4 "54÷9="   # was 62÷4=
Set-CellText $t 9 5 "93÷4="   # was 33÷4=

Set-CellText $t 13 1 "22÷6="  # was 86÷9=
Set-CellText $t 13 2 "85÷8="  # was 27÷6=
Set-CellText $t 13 3 "28÷4="  # was 27÷8=
Set-CellText $t 13 4 "73÷2="  # was 61÷5=
Set-CellText $t 13 5 "97÷5="  # was 32÷9=

Set-CellText $t 17 1 "84÷9="  # was 48÷6=
Set-CellText $t 17 2 "88÷3="  # was 22÷5=
Set-CellText $t 17 3 "86÷5="  # was 24÷9=
Set-CellText $t 17 4 "36÷5="  # was 73÷4=
Set-CellText $t 17 5 "29÷5="  # was 94÷5=

Write-Host "Done updating table cells."
